# "week 4 day 3" — fill in the day-3 VWC_perc readings for both the PIPO
# and PSME tabs (col B), and update each sheet's last-touched
# cursor/scroll/zoom state to match where data entry left off.

$wb = $excel.ActiveWorkbook

# --- Sheet "PIPO" : rows 11-21 (TreeID PIPO11 .. PIPO58) ---
$ws1 = $wb.Worksheets.Item("PIPO")

$ws1.Cells.Item(11, 2).Value = 6.1
$ws1.Cells.Item(12, 2).Value = 10.7
$ws1.Cells.Item(13, 2).Value = 0.6
$ws1.Cells.Item(14, 2).Value = 9.3
$ws1.Cells.Item(15, 2).Value = 8.8
$ws1.Cells.Item(16, 2).Value = 8.2
$ws1.Cells.Item(17, 2).Value = 7.3
$ws1.Cells.Item(18, 2).Value = 8.8
$ws1.Cells.Item(19, 2).Value = 0
$ws1.Cells.Item(20, 2).Value = 10.9
$ws1.Cells.Item(21, 2).Value = 3.2

# --- Sheet "PSME" : rows 2-21 (TreeID PSME3 .. PSME58) ---
$ws2 = $wb.Worksheets.Item("PSME")

$ws2.Cells.Item(2, 2).Value = 11.1
$ws2.Cells.Item(3, 2).Value = 8.4
$ws2.Cells.Item(4, 2).Value = 10.8
$ws2.Cells.Item(5, 2).Value = 0
$ws2.Cells.Item(6, 2).Value = 9.4
$ws2.Cells.Item(7, 2).Value = 2.1
$ws2.Cells.Item(8, 2).Value = 13.3
$ws2.Cells.Item(9, 2).Value = 0
$ws2.Cells.Item(10, 2).Value = 13
$ws2.Cells.Item(11, 2).Value = 1.6
$ws2.Cells.Item(12, 2).Value = 6.6
$ws2.Cells.Item(13, 2).Value = 5.2
$ws2.Cells.Item(14, 2).Value = 8.7
$ws2.Cells.Item(15, 2).Value = 4.7
$ws2.Cells.Item(16, 2).Value = 12.1
$ws2.Cells.Item(17, 2).Value = 8
$ws2.Cells.Item(18, 2).Value = 9.2
$ws2.Cells.Item(19, 2).Value = 0
$ws2.Cells.Item(20, 2).Value = 7.7
$ws2.Cells.Item(21, 2).Value = 2.7

# --- View state: data entry finished on PSME!B21, scrolled & zoomed in,
#     then the user flipped back to the PIPO tab (still scrolled down to
#     where they left off, also resting on B21). ---

$ws2.Activate()
$excel.ActiveWindow.Zoom = 275
$excel.ActiveWindow.ScrollRow = 15
$ws2.Range("B21").Select()

$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 11
$ws1.Range("B21").Select()
